$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Validaciones")

# --- Update task statuses from "Pendiente" to "Terminado" (shared strings 12 -> 11) ---
$ws1.Range("B13").Value = "Terminado"
$ws1.Range("B15").Value = "Terminado"
$ws1.Range("B16").Value = "Terminado"

# --- Row 16 picks up the same highlight formatting already used by row 15 ---
$ws1.Range("A15:B15").Copy()
$ws1.Range("A16:B16").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Rows 15 & 16 are shorter now (single-line wrap instead of tall paragraph) ---
$ws1.Rows(15).RowHeight = 30
$ws1.Rows(16).RowHeight = 30

# --- Hoja1 becomes the active sheet/tab instead of Validaciones ---
$ws1.Activate()
$ws1.Range("A17").Select()
